# Add a new "CatalogueManagement" sheet (first test for catalogue management)
# right after the existing "Admin" sheet, reusing Admin's layout/formatting,
# and make it the active sheet.

$wb = $excel.ActiveWorkbook
$admin = $wb.Worksheets.Item("Admin")

# Duplicate the Admin sheet (keeps its column widths / styles / header rows)
# and place the copy immediately after Admin.
$admin.Copy($null, $admin)

$new = $wb.Worksheets.Item(2)
$new.Name = "CatalogueManagement"

# First catalogue-management test row.
$new.Range("A3").Value = "ViewCatalogueManagementPage"
$new.Range("B3").Value = "positive"
$new.Range("C3").Value = "Nazar_Lelyak"
$new.Range("D3").Value = "qwerty123"

# Restore Admin's own view state (zoom + selection) before switching away.
$admin.Activate()
$excel.ActiveWindow.Zoom = 90
[void]$admin.Range("A3").Select()

# Make the new sheet the active tab.
$new.Activate()
$excel.ActiveWindow.Zoom = 90
[void]$new.Range("F4").Select()
